$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "RF"

# Advance the internal sheetId counter with throwaway add/delete pairs so the
# sheets we actually keep land on sheetId 2..5 (matching the target file,
# which evidently went through some sheet churn before settling).
$t1 = $wb.Worksheets.Add($null, $ws1)
$t2 = $wb.Worksheets.Add($null, $ws1)
$t3 = $wb.Worksheets.Add($null, $ws1)
$t1.Delete()
$t2.Delete()
$t3.Delete()

# Clone RF (while it still holds its original values/styles) for each of the
# other four models, inserting right after RF every time so the final tab
# order becomes RF, ET, AB, GB, XGB.
$ws1.Copy($null, $ws1)
$wsXGB = $wb.Worksheets.Item(2)
$wsXGB.Name = "XGB"

$ws1.Copy($null, $ws1)
$wsGB = $wb.Worksheets.Item(2)
$wsGB.Name = "GB"

$ws1.Copy($null, $ws1)
$wsAB = $wb.Worksheets.Item(2)
$wsAB.Name = "AB"

$ws1.Copy($null, $ws1)
$wsET = $wb.Worksheets.Item(2)
$wsET.Name = "ET"

# ---------------------------------------------------------------------
# RF: update GridSearchCV numbers for the Random Forest model.
# ---------------------------------------------------------------------
$ws1.Range("C4").Value = 11
$ws1.Range("B6").Value = 0.75
$ws1.Range("C6").Value = 0.74
$ws1.Range("E6").Value = 0.99
$ws1.Range("E8").Formula = "=E6/(E6+D4)"
$ws1.Range("E8").NumberFormat = $ws1.Range("D8").NumberFormat
$ws1.Range("D4").Select()

# ---------------------------------------------------------------------
# ET: Extra Trees model.
# ---------------------------------------------------------------------
$wsET.Range("B4").Value = 6
$wsET.Range("C4").Value = 11
$wsET.Range("D4").Formula = "=(B4)/C4"
$wsET.Range("B6").Value = 0.812
$wsET.Range("C6").Value = 0.83
$wsET.Range("E6").Value = 0.94
$wsET.Range("E8").Formula = "=E6/(E6+D4)"
$wsET.Range("E8").NumberFormat = $wsET.Range("D8").NumberFormat
$wsET.Range("B5").Select()

# ---------------------------------------------------------------------
# AB: AdaBoost model.
# ---------------------------------------------------------------------
$wsAB.Range("B4").Value = 5
$wsAB.Range("C4").Value = 12
$wsAB.Range("D4").Formula = "=(C4-B4)/C4"
$wsAB.Range("B6").Value = 0.812
$wsAB.Range("C6").Value = 0.76
$wsAB.Range("E6").Value = 0.99
$wsAB.Range("E8").Formula = "=E6/(E6+D4)"
$wsAB.Range("E8").NumberFormat = $wsAB.Range("D8").NumberFormat
$wsAB.Range("D4").Select()

# ---------------------------------------------------------------------
# GB: Gradient Boosting model.
# ---------------------------------------------------------------------
$wsGB.Range("B4").Value = 7
$wsGB.Range("C4").Value = 12
$wsGB.Range("D4").Formula = "=(C4-B4)/C4"
$wsGB.Range("B6").Value = 0.76
$wsGB.Range("C6").Value = 0.74
$wsGB.Range("E6").Value = 0.99
$wsGB.Range("E8").Formula = "=E6/(E6+D4)"
$wsGB.Range("E8").NumberFormat = $wsGB.Range("D8").NumberFormat
$wsGB.Range("E6:E8").Select()

# ---------------------------------------------------------------------
# XGB: XGBoost model (baseline). D6 stays a plain number, not a formula.
# ---------------------------------------------------------------------
$wsXGB.Range("B4").Value = 8
$wsXGB.Range("C4").Value = 12
$wsXGB.Range("D4").Formula = "=(C4-B4)/C4"
$wsXGB.Range("B6").Value = 0.77
$wsXGB.Range("C6").Value = 0.75
$wsXGB.Range("D6").Value = 0.0099
$wsXGB.Range("E6").Value = 0.99
$wsXGB.Range("E8").Formula = "=E6/(E6+D4)"
$wsXGB.Range("E8").NumberFormat = $wsXGB.Range("D8").NumberFormat
$wsXGB.Range("O10").Select()

# Final active tab is ET, matching the target workbook.
$wsET.Activate()
